$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 3 (Wild Card round results added to Regular season "R" totals) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 530
$wsOff.Range("C3").Value = 367
$wsOff.Range("D3").Value = 108
$wsOff.Range("E3").Value = 42
$wsOff.Range("F3").Value = 7
$wsOff.Range("G3").Value = 8

# --- DEF sheet: update row 3 ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 524
$wsDef.Range("C3").Value = 361
$wsDef.Range("D3").Value = 132
$wsDef.Range("E3").Value = 64
$wsDef.Range("F3").Value = 7
$wsDef.Range("G3").Value = 8
